$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 216 - this shifts the existing rows 216:338 down to 217:339,
# copying formatting from the row above (row 215) by default in Excel.
$ws.Rows("216:216").Insert()

# Copy the fully-formatted row 217 (the old row 216, now shifted down) onto the new row 216
# so that formatting/styles for the new row match the rest of the data block.
$ws.Rows("217:217").Copy()
$ws.Rows("216:216").PasteSpecial()

# Populate the new row 216 with this week's new data point.
$ws.Cells.Item(216, 4).Value = 45097   # D: Fecha
$ws.Cells.Item(216, 10).Value = 43     # J: Volumen
$ws.Cells.Item(216, 11).Value = 19000  # K: Precio minimo
$ws.Cells.Item(216, 12).Value = 20000  # L: Precio maximo
$ws.Cells.Item(216, 13).Value = 19535  # M: Precio promedio ponderado
$ws.Cells.Item(216, 16).Value = 6512   # P: Precio $/Kg
